$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, centered, bordered) from the existing
# last header cell (AC1) onto the three new header cells before setting
# their text, so AD1:AF1 match the look of the other headers.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill every data row (2 through 49) with the constant team record values
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 67   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 95   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
